# The commit removed two MAG rows (even_MAG-GUT47330.fa, originally row 5,
# and even_MAG-GUT54831.fa, originally row 7) from the sheet, shifting the
# remaining rows up so the used range shrinks from A1:AA8 to A1:AA6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the lower row first so the earlier row's index isn't affected by
# the shift caused by the later deletion.
$ws.Rows(7).Delete()
$ws.Rows(5).Delete()
